$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1331.1
$ws.Range("I17").Value = 1750
$ws.Range("J17").Value = 1316.6552
$ws.Range("K17").Value = 5250
$ws.Range("L17").Value = 3949.9656
$ws.Range("M17").Value = -5082
$ws.Range("N17").Value = -4285.9656
$ws.Range("H62").Value = 3992.3635
$ws.Range("I62").Value = 3842.8572
$ws.Range("J62").Value = 4254
$ws.Range("K62").Value = 3842.8572
$ws.Range("L62").Value = 4254
$ws.Range("M62").Value = -3218.8572
$ws.Range("N62").Value = -5502
$ws.Range("H64").Value = 3973.182
$ws.Range("J64").Value = 4440.6
$ws.Range("L64").Value = 4440.6
$ws.Range("N64").Value = -4936.6
$ws.Range("H65").Value = 3992.3635
$ws.Range("I65").Value = 3842.8572
$ws.Range("J65").Value = 4254
$ws.Range("K65").Value = 19214.286
$ws.Range("L65").Value = 21270
$ws.Range("M65").Value = -16094.286
$ws.Range("N65").Value = -27510
$ws.Range("H67").Value = 3973.182
$ws.Range("J67").Value = 4440.6
$ws.Range("L67").Value = 4440.6
$ws.Range("M67").Value = -2725.6667
$ws.Range("N67").Value = -6156.6
$ws.Range("H76").Value = 6375.273
$ws.Range("I76").Value = 4356.4
$ws.Range("J76").Value = 8057.6665
$ws.Range("K76").Value = 4356.4
$ws.Range("L76").Value = 8057.6665
$ws.Range("M76").Value = -4041.4
$ws.Range("N76").Value = -8687.666499999999
$ws.Range("H79").Value = 6375.273
$ws.Range("I79").Value = 4356.4
$ws.Range("J79").Value = 8057.6665
$ws.Range("K79").Value = 4356.4
$ws.Range("L79").Value = 8057.6665
$ws.Range("M79").Value = -3264.4
$ws.Range("N79").Value = -10241.6665
$ws.Range("H111").Value = 1246.5714
$ws.Range("I111").Value = 937.9167
$ws.Range("J111").Value = 1658.1111
$ws.Range("K111").Value = 2813.7501
$ws.Range("L111").Value = 4974.3333
$ws.Range("M111").Value = 253.2498999999998
$ws.Range("N111").Value = -11108.3333
$ws.Range("H113").Value = 2952.3
$ws.Range("I113").Value = 2996.3333
$ws.Range("J113").Value = 2933.4285
$ws.Range("K113").Value = 2996.3333
$ws.Range("L113").Value = 2933.4285
$ws.Range("M113").Value = 257.6667000000002
$ws.Range("N113").Value = -9441.4285
$ws.Range("H116").Value = 7272.5835
$ws.Range("J116").Value = 5301.8
$ws.Range("L116").Value = 5301.8
$ws.Range("N116").Value = -12185.8
$ws.Range("H132").Value = 4127.5205
$ws.Range("I132").Value = 2074.6667
$ws.Range("K132").Value = 6224.000100000001
$ws.Range("M132").Value = -3694.000100000001
$ws.Range("H135").Value = 1394.5555
$ws.Range("I135").Value = 1035.4117
$ws.Range("K135").Value = 9318.705300000001
$ws.Range("M135").Value = -6783.705300000001
$ws.Range("H137").Value = 4652.75
$ws.Range("I137").Value = 4771.3706
$ws.Range("J137").Value = 1450
$ws.Range("K137").Value = 14314.1118
$ws.Range("L137").Value = 4350
$ws.Range("M137").Value = -11764.1118
$ws.Range("N137").Value = -9450
$ws.Range("H138").Value = 2455.4714
$ws.Range("I138").Value = 2129.9614
$ws.Range("J138").Value = 2647.818
$ws.Range("K138").Value = 6389.8842
$ws.Range("L138").Value = 7943.454000000001
$ws.Range("M138").Value = -1249.8842
$ws.Range("N138").Value = -18223.454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5740.923
$ws.Range("I61").Value = 9924.286
$ws.Range("J61").Value = 4199.684
$ws.Range("K61").Value = 9924.286
$ws.Range("L61").Value = 4199.684
$ws.Range("M61").Value = -9712.286
$ws.Range("N61").Value = -4623.684
$ws.Range("H63").Value = 1969.4
$ws.Range("I63").Value = 1887.8889
$ws.Range("J63").Value = 2703
$ws.Range("K63").Value = 1887.8889
$ws.Range("L63").Value = 2703
$ws.Range("M63").Value = -1201.8889
$ws.Range("N63").Value = -4075
$ws.Range("H66").Value = 1969.4
$ws.Range("I66").Value = 1887.8889
$ws.Range("J66").Value = 2703
$ws.Range("K66").Value = 9439.4445
$ws.Range("L66").Value = 13515
$ws.Range("M66").Value = -6007.4445
$ws.Range("N66").Value = -20379
$ws.Range("H122").Value = 1767.1818
$ws.Range("I122").Value = 1868.3448
$ws.Range("J122").Value = 1571.6
$ws.Range("K122").Value = 5605.0344
$ws.Range("L122").Value = 4714.799999999999
$ws.Range("M122").Value = -3155.0344
$ws.Range("N122").Value = -9614.799999999999
$ws.Range("H132").Value = 5954.7
$ws.Range("I132").Value = 5364.143
$ws.Range("J132").Value = 7332.6665
$ws.Range("K132").Value = 16092.429
$ws.Range("L132").Value = 21997.9995
$ws.Range("M132").Value = -13562.429
$ws.Range("N132").Value = -27057.9995
$ws.Range("H136").Value = 5740.923
$ws.Range("I136").Value = 9924.286
$ws.Range("J136").Value = 4199.684
$ws.Range("K136").Value = 29772.858
$ws.Range("L136").Value = 12599.052
$ws.Range("M136").Value = -27222.858
$ws.Range("N136").Value = -17699.052

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1787.1621
$ws.Range("I105").Value = 1768.7858
$ws.Range("J105").Value = 1844.3334
$ws.Range("K105").Value = 1768.7858
$ws.Range("L105").Value = 1844.3334
$ws.Range("M105").Value = -21.78580000000011
$ws.Range("N105").Value = -5338.3334
$ws.Range("H107").Value = 2699.0605
$ws.Range("I107").Value = 2322.9614
$ws.Range("J107").Value = 4096
$ws.Range("K107").Value = 2322.9614
$ws.Range("L107").Value = 4096
$ws.Range("M107").Value = -402.9614000000001
$ws.Range("N107").Value = -7936
$ws.Range("H134").Value = 3363
$ws.Range("I134").Value = 3668.125
$ws.Range("J134").Value = 3188.6428
$ws.Range("K134").Value = 11004.375
$ws.Range("L134").Value = 9565.928400000001
$ws.Range("M134").Value = -8469.375
$ws.Range("N134").Value = -14635.9284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7064.1763
$ws.Range("I22").Value = 9902
$ws.Range("K22").Value = 9902
$ws.Range("M22").Value = -9552
$ws.Range("H62").Value = 4970.3477
$ws.Range("I62").Value = 4627.6665
$ws.Range("J62").Value = 5344.1816
$ws.Range("K62").Value = 4627.6665
$ws.Range("L62").Value = 5344.1816
$ws.Range("M62").Value = -4003.6665
$ws.Range("N62").Value = -6592.1816
$ws.Range("H65").Value = 4970.3477
$ws.Range("I65").Value = 4627.6665
$ws.Range("J65").Value = 5344.1816
$ws.Range("K65").Value = 23138.3325
$ws.Range("L65").Value = 26720.908
$ws.Range("M65").Value = -20018.3325
$ws.Range("N65").Value = -32960.908
$ws.Range("H107").Value = 859.65625
$ws.Range("I107").Value = 716.7143
$ws.Range("K107").Value = 716.7143
$ws.Range("M107").Value = 1203.2857
$ws.Range("H131").Value = 83416.46000000001
$ws.Range("J131").Value = 83701.164
$ws.Range("L131").Value = 83701.164
$ws.Range("N131").Value = -93781.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2667.4666
$ws.Range("I11").Value = 581.5
$ws.Range("J11").Value = 3426
$ws.Range("K11").Value = 1744.5
$ws.Range("L11").Value = 10278
$ws.Range("M11").Value = -1604.5
$ws.Range("N11").Value = -10558
$ws.Range("H94").Value = 24232988
$ws.Range("I94").Value = 1999
$ws.Range("K94").Value = 5997
$ws.Range("M94").Value = -5321
$ws.Range("H132").Value = 2666.5
$ws.Range("J132").Value = 3547.3
$ws.Range("L132").Value = 31925.7
$ws.Range("N132").Value = -36985.7
$ws.Range("H141").Value = 7030
$ws.Range("I141").Value = 7030
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 21090
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -15910
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 21424.416
$ws.Range("I70").Value = 41816.2
$ws.Range("K70").Value = 41816.2
$ws.Range("M70").Value = -41546.2
$ws.Range("H73").Value = 21424.416
$ws.Range("I73").Value = 41816.2
$ws.Range("K73").Value = 41816.2
$ws.Range("M73").Value = -40880.2
$ws.Range("H80").Value = 53337164
$ws.Range("I80").Value = 80003220
$ws.Range("J80").Value = 5066.6665
$ws.Range("K80").Value = 80003220
$ws.Range("L80").Value = 5066.6665
$ws.Range("M80").Value = -80002222
$ws.Range("N80").Value = -7062.6665
$ws.Range("H83").Value = 53337164
$ws.Range("I83").Value = 80003220
$ws.Range("J83").Value = 5066.6665
$ws.Range("K83").Value = 400016100
$ws.Range("L83").Value = 25333.3325
$ws.Range("M83").Value = -400011108
$ws.Range("N83").Value = -35317.3325
$ws.Range("H132").Value = 6966
$ws.Range("I132").Value = 9999
$ws.Range("J132").Value = 5449.5
$ws.Range("K132").Value = 29997
$ws.Range("L132").Value = 16348.5
$ws.Range("M132").Value = -27467
$ws.Range("N132").Value = -21408.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7209560
$ws.Range("I136").Value = 9010788
$ws.Range("J136").Value = 4650
$ws.Range("K136").Value = 27032364
$ws.Range("L136").Value = 13950
$ws.Range("M136").Value = -27029814
$ws.Range("N136").Value = -19050

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11589
$ws.Range("J45").Value = 11558.25
$ws.Range("L45").Value = 11558.25
$ws.Range("N45").Value = -12540.25
$ws.Range("H124").Value = 69348.836
$ws.Range("J124").Value = 69348.836
$ws.Range("L124").Value = 69348.836
$ws.Range("N124").Value = -79168.836
$ws.Range("H132").Value = 13999.833
$ws.Range("I132").Value = 14799.8
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 44399.39999999999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -41869.39999999999
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 2828.1853
$ws.Range("I136").Value = 3495.9285
$ws.Range("J136").Value = 2109.077
$ws.Range("K136").Value = 10487.7855
$ws.Range("L136").Value = 6327.231000000001
$ws.Range("M136").Value = -7937.7855
$ws.Range("N136").Value = -11427.231
